$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.906.73'
$ws.Range("E2").Value = '  -1.02%  '

$ws.Range("D3").Value = '2.337.83'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.51'
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.35'
$ws.Range("E6").Value = '  -3.70%  '

$ws.Range("E7").Value = '  -1.13%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.09'
$ws.Range("E10").Value = '  -4.36%  '

$ws.Range("E11").Value = '  -1.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.66'
$ws.Range("E12").Value = '  -4.89%  '

$ws.Range("E13").Value = '  +1.58%  '

$ws.Range("E14").Value = '  -2.57%  '

$ws.Range("D15").Value = '2.701.68'
$ws.Range("E15").Value = '  +0.42%  '

$ws.Range("D16").Value = '2.347.85'
$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").Value = '42.842.00'
$ws.Range("E18").Value = '  -0.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.06'
$ws.Range("E19").Value = '  -5.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  +1.54%  '

$ws.Range("D21").Value = '0.0₃0888'
$ws.Range("E21").Value = '  -1.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.84'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.24'
$ws.Range("E23").Value = '  -0.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  -1.92%  '

$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("E26").Value = '  -1.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.53'
$ws.Range("E27").Value = '  -2.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  +7.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.13'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.35'
$ws.Range("E30").Value = '  -5.64%  '

$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("E32").Value = '  -0.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0733'
$ws.Range("E33").Value = '  +3.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.23'
$ws.Range("E34").Value = '  -3.53%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.33'
$ws.Range("E35").Value = '  -1.06%  '

$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.36'
$ws.Range("E36").Value = '  -3.24%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("E37").Value = '  +2.22%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '125.22'
$ws.Range("E38").Value = '  -23.75%  '

$ws.Range("E40").Value = '  -1.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.21'
$ws.Range("E41").Value = '  +15.14%  '

$ws.Range("E42").Value = '  -1.52%  '

$ws.Range("D43").Value = '1.937.90'
$ws.Range("E43").Value = '  -2.46%  '

$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.14'
$ws.Range("E45").Value = '  -4.51%  '

$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("E47").Value = '  -3.47%  '

$ws.Range("D48").Value = '2.566.12'
$ws.Range("E48").Value = '  +0.40%  '

$ws.Range("E49").Value = '  +0.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.75'
$ws.Range("E50").Value = '  -2.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.64'
$ws.Range("E51").Value = '  -1.54%  '
